# test carga masiva eventos
# Replace the sample rows (2 and 3) on the events template with new test data,
# clear the extra height on row 3, and move the active selection to I14.
# NOTE: the cell write order below mirrors how the data was originally typed in
# (interleaved across rows 2/3), which determines the resulting shared-string order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Clase 1"
$ws.Range("A3").Value = "Clase 1"

$ws.Range("B2").Value = "empresarial"
$ws.Range("B3").Value = "deportivo"

$ws.Range("C2").Value = "presencial"
$ws.Range("C3").Value = "virtual"

$ws.Range("D2").Value = "clase 1 1 1 1 1"
$ws.Range("D3").Value = "clase 2 2 2 2 2"

$ws.Range("E2").Value = "Barranquilla"
$ws.Range("E3").Value = "no aplica"

$ws.Range("F3").Value = "url//"
$ws.Range("F2").Value = "dg 59 # 45-65"

$ws.Range("G2").Value = 3118522123
$ws.Range("G3").Value = 3118522123

$ws.Range("H2").Value = "No aplican requisitos"
$ws.Range("H3").Value = "internet y equipo"

$ws.Range("I2").Value = 10
$ws.Range("I3").Value = 40

# Row 3 no longer needs the taller wrapped-text height
$ws.Rows(3).EntireRow.AutoFit()

# Update the selected cell shown in the sheet view
$ws.Range("I14").Select()
